$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = '28.414.71'
$ws.Range("E2").Value = '  +0.40%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = '1.828.12'
$ws.Range("E3").Value = '  +2.13%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.28'
$ws.Range("E5").Value = '  +0.56%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range("E6").Value = '  -0.04%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5320'
$ws.Range("E7").Value = '  -0.68%  '

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4046'
$ws.Range("E8").Value = '  +7.62%  '

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07613'
$ws.Range("E9").Value = '  +1.85%  '

# Row 10: 'OKB' -> 'OKB'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.82'
$ws.Range("E10").Value = '  +1.24%  '

# Row 11: 'Polygon' -> 'Polygon'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.106'
$ws.Range("E11").Value = '  +1.32%  '

# Row 12: 'Polkadot' -> 'Polkadot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.351'
$ws.Range("E12").Value = '  +4.30%  '

# Row 14: 'Solana' -> 'Solana'
$ws.Range("E14").Value = '  +2.20%  '

# Row 15: 'Chainlink' -> 'Chainlink'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.535'
$ws.Range("E15").Value = '  +4.06%  '

# Row 16: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D16").Value = '1.825.88'
$ws.Range("E16").Value = '  +2.25%  '

# Row 17: 'Litecoin' -> 'ShibaInu'
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001074'
$ws.Range("E17").Value = '  +1.68%  '

# Row 18: 'ShibaInu' -> 'Litecoin'
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.27'
$ws.Range("E18").Value = '  +0.20%  '

# Row 19: 'TRON' -> 'TRON'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06609'
$ws.Range("E19").Value = '  +1.82%  '

# Row 20: 'Avalanche' -> 'Avalanche'
$ws.Range("E20").Value = '  +1.36%  '

# Row 21: 'Dai' -> 'Dai'
$ws.Range("E21").Value = '  -0.02%  '

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.068'
$ws.Range("E22").Value = '  +2.49%  '

# Row 23: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D23").Value = '28.451.55'
$ws.Range("E23").Value = '  +0.45%  '

# Row 24: 'Cosmos' -> 'Cosmos'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.28'
$ws.Range("E24").Value = '  +1.84%  '

# Row 25: 'Toncoin' -> 'Toncoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.157'
$ws.Range("E25").Value = '  +3.11%  '

# Row 26: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.479'
$ws.Range("E26").Value = '  +8.36%  '

# Row 27: 'Monero' -> 'Monero'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.82'
$ws.Range("E27").Value = '  -0.95%  '

# Row 28: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("E28").Value = '  +1.59%  '

# Row 29: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D29").Value = '2.030.40'
$ws.Range("E29").Value = '  +1.93%  '

# Row 30: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.56'
$ws.Range("E30").Value = '  +1.60%  '

# Row 31: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.123'
$ws.Range("E31").Value = '  +2.82%  '

# Row 32: 'Stellar' -> 'Stellar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1093'
$ws.Range("E32").Value = '  +4.48%  '

# Row 33: 'Filecoin' -> 'Filecoin'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.679'
$ws.Range("E33").Value = '  +2.94%  '

# Row 34: 'HuobiToken' -> 'HuobiToken'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.658'
$ws.Range("E34").Value = '  -0.11%  '

# Row 35: 'Hedera' -> 'Hedera'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07167'
$ws.Range("E35").Value = '  +11.70%  '

# Row 36: 'Algorand' -> 'Algorand'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2262'
$ws.Range("E36").Value = '  +0.49%  '

# Row 37: 'VeChain' -> 'VeChain'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02342'
$ws.Range("E37").Value = '  +3.03%  '

# Row 38: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.229'
$ws.Range("E38").Value = '  +4.58%  '

# Row 39: 'FraxShare' -> 'FraxShare'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.772'
$ws.Range("E39").Value = '  +3.36%  '

# Row 40: 'TheSandbox' -> 'TheSandbox'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6274'
$ws.Range("E40").Value = '  +2.08%  '

# Row 41: 'Aptos' -> 'Aptos'
$ws.Range("E41").Value = '  +2.92%  '

# Row 42: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.184'
$ws.Range("E42").Value = '  +0.52%  '

# Row 43: 'Frax' -> 'Frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.07%  '

# Row 44: 'WEMIXTOKEN' -> 'WEMIXTOKEN'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.399'
$ws.Range("E44").Value = '  -2.62%  '

# Row 45: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.42'
$ws.Range("E45").Value = '  +0.68%  '

# Row 46: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.705'
$ws.Range("E46").Value = '  +1.14%  '

# Row 47: 'Decentraland' -> 'Decentraland'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5847'
$ws.Range("E47").Value = '  +1.57%  '

# Row 48: 'Quant' -> 'Quant'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.11'
$ws.Range("E48").Value = '  +0.46%  '

# Row 49: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.989'
$ws.Range("E49").Value = '  +3.14%  '

# Row 50: 'EOS' -> 'EOS'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.197'
$ws.Range("E50").Value = '  +0.01%  '

# Row 51: 'Cronos' -> 'Cronos'
$ws.Range("E51").Value = '  +0.82%  '
